# PaperChoices.xlsx update:
#  - PaperList: add "Analysis" (E) and "PDF(alternate link)" (F) columns with per-study data
#  - Assigned: re-sort the assignment table by PaperNum (ascending)
#  - Selections / active sheet updated to match the authoring session

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PaperList")
$ws2 = $wb.Worksheets.Item("Group_Members")
$ws4 = $wb.Worksheets.Item("Assigned")

# PaperList already carried left-over column-width formatting on (empty) columns E:G.
# Insert a fresh column at F so that new column F becomes a brand new entry and the
# old left-over E/F/G formatting lands on E (reused in place) / G / H, matching the
# target column layout.
$ws1.Columns.Item(6).EntireColumn.Insert()

# --- PaperList: fill column E ("Analysis") for the data rows first, then the header ---
$ws1.Range("E2").Value = "one-way-anova"
$ws1.Range("E3").Value = "one-sample-t-test,dependent-samples-t-test,multiple-linear-regression"
$ws1.Range("E4").Value = "multiple-linear-regression"
$ws1.Range("E5").Value = "independent-samples-t-test,multiple-linear-regression,simple-linear-regression,dependent-samples-t-test"
$ws1.Range("E6").Value = "dependent-samples-t-test,one-sample-t-test"
$ws1.Range("E7").Value = "multiple-linear-regression"
$ws1.Range("E8").Value = "two-way-anova,one-way-anova"
$ws1.Range("E9").Value = "independent-samples-t-test"
$ws1.Range("E10").Value = "dependent-samples-t-test"
$ws1.Range("E11").Value = "dependent-samples-t-test"
$ws1.Range("E12").Value = "independent-samples-t-test"
$ws1.Range("E13").Value = "one-way-anova"
$ws1.Range("E14").Value = "multiple-linear-regression"
$ws1.Range("E15").Value = "one-way-anova"
$ws1.Range("E16").Value = "simple-linear-regression,two-way-anova,independent-samples-t-test"
$ws1.Range("E1").Value = "Analysis"

# --- PaperList: fill column F ("PDF(alternate link)") for the data rows first, then the header ---
$ws1.Range("F2").Value = "0956797614545886.pdf (https://dl.airtable.com/sQiEblcHTzqA60UfEhKb_0956797614545886.pdf)"
$ws1.Range("F3").Value = "0956797615626691.pdf (https://dl.airtable.com/AKd7LiGUQVmsluL6mPKN_0956797615626691.pdf)"
$ws1.Range("F4").Value = "journal.pone.0182159.pdf (https://dl.airtable.com/9dKhg8EQnuaGoAtKob0d_journal.pone.0182159.pdf)"
$ws1.Range("F5").Value = "journal.pone.0182817.pdf (https://dl.airtable.com/Xer6TyfTS1yudPIFe5IL_journal.pone.0182817.pdf)"
$ws1.Range("F6").Value = "0956797614533801.pdf (https://dl.airtable.com/osAdTJQTUmMZjRkvevOk_0956797614533801.pdf)"
$ws1.Range("F7").Value = "journal.pone.0182239.pdf (https://dl.airtable.com/nAle73A2Rpy1JM4CMsQm_journal.pone.0182239.pdf)"
$ws1.Range("F8").Value = "journal.pone.0182907.pdf (https://dl.airtable.com/IZvp6MTPiGSoFZ8oeElA_journal.pone.0182907.pdf)"
$ws1.Range("F9").Value = "0956797615620784.pdf (https://dl.airtable.com/tfI9lhxORP2jJ1A0xwO2_0956797615620784.pdf)"
$ws1.Range("F12").Value = "journal.pone.0173493.pdf (https://dl.airtable.com/6N3azAFzROGVbLL6VvSB_journal.pone.0173493.pdf)"
$ws1.Range("F13").Value = "journal.pone.0177758.pdf (https://dl.airtable.com/9QV0nDO6SXqqzRKx18fg_journal.pone.0177758.pdf)"
$ws1.Range("F15").Value = "journal.pone.0152576.PDF (https://dl.airtable.com/P2I3CFfNTDmNV7nxo53h_journal.pone.0152576.PDF)"
$ws1.Range("F1").Value = "PDF(alternate link)"

# Auto-size the two new columns to fit their (much wider) content
$ws1.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws1.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# --- Assigned: re-sort the whole table by column C (PaperNum), ascending ---
$sortRange = $ws4.Range("A2:C10")
$sortKey   = $ws4.Range("C2:C10")
$ws4.Sort.SortFields.Clear()
$ws4.Sort.SortFields.Add($sortKey)
$ws4.Sort.SetRange($sortRange)
$ws4.Sort.Header = 0
$ws4.Sort.Apply()

# --- Selections, matching the end-of-session UI state ---
$ws2.Range("B10").Select()
$ws4.Range("C10").Select()
$ws1.Range("A5").Select()
$ws1.Activate()
